# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values. Update rows 2-25 with the regenerated values
# (row 5 is unchanged since its K value was already 0).
$kValues = @{
    2  = 4
    3  = 11
    4  = 8
    6  = 9
    7  = 7
    8  = 2
    9  = 8
    10 = 10
    11 = 5
    12 = 9
    13 = 3
    14 = 3
    15 = 6
    16 = 2
    17 = 7
    18 = 11
    19 = 6
    20 = 6
    21 = 5
    22 = 3
    23 = 3
    24 = 2
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
